$d = $word.ActiveDocument

# Remove paragraphs 10-15 (empty line, Average/Best/Worst/Last-5-games block,
# and the trailing empty line) that sat between the Expectimax stats and the
# Heatmap heading. Delete from the bottom up so earlier indices stay valid.
$start = $d.Paragraphs.Item(10).Range.Start
$end = $d.Paragraphs.Item(15).Range.End
$d.Range($start, $end).Delete()

# Remove the empty paragraph between the Monte Carlo stats and the
# Expectimax heading.
$d.Paragraphs.Item(5).Range.Delete()

# Update the Expectimax statistics.
$d.Content.Find.Execute("Average shots to win: 49.6", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Average shots to win: 49.1", 2)
$d.Content.Find.Execute("Best game (fewest shots): 31", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Best game (fewest shots): 28", 2)
$d.Content.Find.Execute("Worst game (most shots): 77", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Worst game (most shots): 78", 2)

# Update the Heatmap statistics.
$d.Content.Find.Execute("Average shots to win: 51.4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Average shots to win: 48.5", 2)
$d.Content.Find.Execute("Best game (fewest shots): 34", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Best game (fewest shots): 29", 2)
$d.Content.Find.Execute("Worst game (most shots): 79", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Worst game (most shots): 92", 2)
